$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at the top of this block (60-62), shifting old rows 60-81 down to 63-84
$ws.Rows("60:62").Insert()

$ws.Range("A60").Value = 5
$ws.Range("B60").Value = 'Macroferia Regional de Talca'
$ws.Range("C60").Value = 'Maule'
$ws.Range("D60").Value = 44518
$ws.Range("E60").Value = 7
$ws.Range("F60").Value = 'Fruta'
$ws.Range("G60").Value = 100103
$ws.Range("H60").Value = 'Frutos de hueso (carozo)'
$ws.Range("I60").Value = 100103001
$ws.Range("J60").Value = 'Cereza'
$ws.Range("K60").Value = 'Royal Dawn'
$ws.Range("L60").Value = 'Primera'
$ws.Range("M60").Value = 100
$ws.Range("N60").Value = 20000
$ws.Range("O60").Value = 20000
$ws.Range("P60").Value = 20000
$ws.Range("Q60").Value = '$/bandeja 10 kilos'
$ws.Range("R60").Value = 'Provincia de Curicó'
$ws.Range("S60").Value = 2000
$ws.Range("T60").Value = 10

$ws.Range("A61").Value = 5
$ws.Range("B61").Value = 'Macroferia Regional de Talca'
$ws.Range("C61").Value = 'Maule'
$ws.Range("D61").Value = 44518
$ws.Range("E61").Value = 7
$ws.Range("F61").Value = 'Fruta'
$ws.Range("G61").Value = 100103
$ws.Range("H61").Value = 'Frutos de hueso (carozo)'
$ws.Range("I61").Value = 100103001
$ws.Range("J61").Value = 'Cereza'
$ws.Range("K61").Value = 'Royal Dawn'
$ws.Range("L61").Value = 'Segunda'
$ws.Range("M61").Value = 30
$ws.Range("N61").Value = 28000
$ws.Range("O61").Value = 28000
$ws.Range("P61").Value = 28000
$ws.Range("Q61").Value = '$/caja 15 kilos'
$ws.Range("R61").Value = 'Provincia de Curicó'
$ws.Range("S61").Value = 1867
$ws.Range("T61").Value = 15

$ws.Range("A62").Value = 5
$ws.Range("B62").Value = 'Macroferia Regional de Talca'
$ws.Range("C62").Value = 'Maule'
$ws.Range("D62").Value = 44518
$ws.Range("E62").Value = 7
$ws.Range("F62").Value = 'Fruta'
$ws.Range("G62").Value = 100103
$ws.Range("H62").Value = 'Frutos de hueso (carozo)'
$ws.Range("I62").Value = 100103001
$ws.Range("J62").Value = 'Cereza'
$ws.Range("K62").Value = 'Santina'
$ws.Range("L62").Value = 'Segunda'
$ws.Range("M62").Value = 150
$ws.Range("N62").Value = 15000
$ws.Range("O62").Value = 18000
$ws.Range("P62").Value = 17000
$ws.Range("Q62").Value = '$/bandeja 10 kilos'
$ws.Range("R62").Value = 'Provincia de Curicó'
$ws.Range("S62").Value = 1700
$ws.Range("T62").Value = 10
